# ---------------------------------------------------------------------------
# Applies the "ubo and directory done" change set:
#   - tweak selections/cursor position on "personal data" and "Business Address"
#   - set a (whitespace) value in "Business Address"!E4
#   - add two new worksheets "UBOs" and "Directors" with validation-test data
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "personal data" sheet: move the view / selection to J5 (scrolled to col H)
# ---------------------------------------------------------------------------
$wsPersonal = $wb.Worksheets.Item("personal data")
$wsPersonal.Activate()
$excel.ActiveWindow.ScrollColumn = 8
$wsPersonal.Range("J5").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. "Business Address" sheet: selection -> D3, E4 gets a blank/space value
# ---------------------------------------------------------------------------
$wsBusiness = $wb.Worksheets.Item("Business Address")
$wsBusiness.Range("E4").Value = "   "
$wsBusiness.Activate()
$wsBusiness.Range("D3").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Add "UBOs" worksheet (after "Business Address")
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsUbo = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsUbo.Name = "UBOs"

$wsUbo.Range("A1").Value = "countryName"
$wsUbo.Range("B1").Value = "uboIdentityNumber"
$wsUbo.Range("C1").Value = "uboAddress"
$wsUbo.Range("D1").Value = "uboPercentage"
$wsUbo.Range("E1").Value = "countyValidation"
$wsUbo.Range("F1").Value = "uboNoValidation"
$wsUbo.Range("G1").Value = "uboAddressValidation"
$wsUbo.Range("H1").Value = "uboPercentageValidation"
$wsUbo.Range("I1").Value = "directorsValidation"

$wsUbo.Range("B2").Value = 123456789
$wsUbo.Range("C2").Value = "22, Choyangmen"
$wsUbo.Range("D2").Value = 46
$wsUbo.Range("E2").Value = "Nationality is required."

$wsUbo.Range("A3").Value = "Indian"
$wsUbo.Range("C3").Value = "22, Choyangmen"
$wsUbo.Range("D3").Value = 26
$wsUbo.Range("F3").Value = "Identity number is required."

$wsUbo.Range("A4").Value = "Indian"
$wsUbo.Range("B4").Value = 123
$wsUbo.Range("C4").Value = "22, Choyangmen"
$wsUbo.Range("D4").Value = 26
$wsUbo.Range("F4").Value = "Identity number should be at least 5 characters long."

$wsUbo.Range("A5").Value = "Indian"
$wsUbo.Range("B5").Value = 123456
$wsUbo.Range("D5").Value = 26
$wsUbo.Range("G5").Value = "Full address is required."

$wsUbo.Range("A6").Value = "Indian"
$wsUbo.Range("B6").Value = 123456
$wsUbo.Range("C6").Value = "22 test"
$wsUbo.Range("D6").Value = 26
$wsUbo.Range("G6").Value = "Full address should be at least 10 characters long."

$wsUbo.Range("A7").Value = "Indian"
$wsUbo.Range("B7").Value = 123456
$wsUbo.Range("C7").Value = "test22, Choyangmen"
$wsUbo.Range("H7").Value = "Percentage of ownership is required."

$wsUbo.Range("A8").Value = "Indian"
$wsUbo.Range("B8").Value = 123456
$wsUbo.Range("C8").Value = "test22, Choyangmen"
$wsUbo.Range("D8").Value = 15
$wsUbo.Range("H8").Value = "Percentage of ownership must be 25 or greater."

$wsUbo.Range("A9").Value = "Indian"
$wsUbo.Range("B9").Value = 123456
$wsUbo.Range("C9").Value = "test22, Choyangmen"
$wsUbo.Range("D9").Value = 25
$wsUbo.Range("I9").Value = "Directors are required."

$wsUbo.Columns.Item(1).ColumnWidth = 22.26
$wsUbo.Columns.Item(2).ColumnWidth = 26.87
$wsUbo.Columns.Item(3).ColumnWidth = 22.76
$wsUbo.Columns.Item(4).ColumnWidth = 30.96
$wsUbo.Columns.Item(5).ColumnWidth = 20.07
$wsUbo.Columns.Item(6).ColumnWidth = 52.96
$wsUbo.Columns.Item(7).ColumnWidth = 48.57
$wsUbo.Columns.Item(8).ColumnWidth = 49.57
$wsUbo.Columns.Item(9).ColumnWidth = 23.76

$wsUbo.Range("G5:G6").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Add "Directors" worksheet (after "UBOs")
# ---------------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsDir = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$wsDir.Name = "Directors"

$wsDir.Range("A1").Value = "countryName"
$wsDir.Range("B1").Value = "directorIdentityNumber"
$wsDir.Range("C1").Value = "addressOfResidence"
$wsDir.Range("D1").Value = "countyValidation"
$wsDir.Range("E1").Value = "directorINoValidation"
$wsDir.Range("F1").Value = "directorAddressValidation"

$wsDir.Range("B2").Value = 123456789
$wsDir.Range("C2").Value = "22, Choyangmen"
$wsDir.Range("D2").Value = "Nationality is required."

$wsDir.Range("A3").Value = "Indian"
$wsDir.Range("C3").Value = "22, Choyangmen"
$wsDir.Range("E3").Value = "Identity number is required."

$wsDir.Range("A4").Value = "Indian"
$wsDir.Range("B4").Value = 123
$wsDir.Range("C4").Value = "22, Choyangmen"
$wsDir.Range("E4").Value = "Identity number should be at least 5 characters long."

$wsDir.Range("A5").Value = "Indian"
$wsDir.Range("B5").Value = 123456789
$wsDir.Range("F5").Value = "Full address is required."

$wsDir.Range("A6").Value = "Indian"
$wsDir.Range("B6").Value = 123456789
$wsDir.Range("C6").Value = 22
$wsDir.Range("F6").Value = "Full address should be at least 10 characters long."

$wsDir.Range("A7").Value = "Indian"
$wsDir.Range("B7").Value = 123456789
$wsDir.Range("C7").Value = "22, Choyangmen test"

$wsDir.Columns.Item(1).ColumnWidth = 34.07
$wsDir.Columns.Item(2).ColumnWidth = 33.57
$wsDir.Columns.Item(3).ColumnWidth = 22.87
$wsDir.Columns.Item(4).ColumnWidth = 34.07
$wsDir.Columns.Item(5).ColumnWidth = 48.57
$wsDir.Columns.Item(6).ColumnWidth = 44.76

$wsDir.Range("A7,D7").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5. Re-activate "personal data" so the saved workbook opens on the original tab
# ---------------------------------------------------------------------------
$wsPersonal.Activate()
